$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.981.03'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.641.04'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('D12').Value = '1.874.14'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.643.03'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('E15').Value = '  +4.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.89'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '27.975.98'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '232.65'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.60'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.75'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.08'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '151.59'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.99'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.72'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0484'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.10'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.412.00'
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.886'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.91%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0169'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.914'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.81%  '
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  +7.17%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '66.32'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.02%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.49'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.12%  '
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').Value = '1.782.61'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '88.03'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.100'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.64'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.82%  '
